$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LH_WF-ADMINHOME-REVIEWS-SHEET")
$ws2 = $wb.Worksheets.Item("VERSION-HISTORY")

# Reviewer verification for LH_REVIEW_WF_ADMINHOME_001 is now closed
$ws1.Range("I2").Value = "closed"

# Record new version history entry: v1.1 - close reviewer varifications
$ws2.Range("A2:D2").Copy()
$ws2.Range("A3:D3").PasteSpecial(-4122)
$ws2.Range("A3").Value = "v1.1"
$ws2.Range("B3").Value = "Hala Eldaly"
$ws2.Range("C3").Value = "close reviewer varifications"
$ws2.Range("D3").Value = $ws2.Range("D2").Value2

$ws2.Range("C4").Select() | Out-Null
$ws1.Range("I14").Select() | Out-Null
